$d = $word.ActiveDocument

# 1. Change the job title from ".NET Software Engineer" to "Software Engineer"
$d.Content.Find.Execute(".NET Software Engineer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Software Engineer", 2)

# 2. Insert " Vue," into the summary sentence after "React," and before " TypeScript",
#    splitting the single run into three runs (matching how Word naturally creates
#    a new run boundary around freshly-typed/inserted text).
$findRng = $d.Content
$findRng.Find.Execute("React,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)                 # wdCollapseEnd -> collapse to right after "React,"
$insertStart = $findRng.Start
$findRng.InsertAfter(" Vue,")
$insertEnd = $insertStart + 5        # length of " Vue,"

# Nudge formatting on just the inserted text (set then restore Bold) so the engine
# keeps it as a distinct run instead of silently re-merging it with its neighbors.
$newRunRng = $d.Range($insertStart, $insertEnd)
$newRunRng.Bold = 1
$newRunRng.Bold = 0
